$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add RSI values (column E) for rows 2-5
$ws.Range("E2").Value = 30.7
$ws.Range("E3").Value = 22.2
$ws.Range("E4").Value = 57.4
$ws.Range("E5").Value = 44.8

# Update 점수(룰) (column G) for rows 4-5
$ws.Range("G4").Value = 30
$ws.Range("G5").Value = 20

# Update 최종점수 (column K) for rows 4-5
$ws.Range("K4").Value = 54.8
$ws.Range("K5").Value = 54.2

# Update MACRO_SCORE (column N) for rows 2-5
$ws.Range("N2").Value = 85.87127175646313
$ws.Range("N3").Value = 85.87127175646313
$ws.Range("N4").Value = 85.87127175646313
$ws.Range("N5").Value = 85.87127175646313
